# Updated Phenology DAS in observed files
# Fill in Cotton.Phenology.SquaringDAS (R), FloweringDAS (S),
# MaturityDAS (V) and HarvestRipeDAS (W) for every observed data row
# (rows 2-44) with their computed DAS values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

for ($row = 2; $row -le 44; $row++) {
    $ws.Range("R$row").Value = 53
    $ws.Range("S$row").Value = 77
    $ws.Range("V$row").Value = 148
    $ws.Range("W$row").Value = 183
}
